# Update countries & provincias Spain
# - Swap the displayed country names for rows 128/129 (Lituania <-> Eslovenia)
#   and rows 213/214 (Islas Malvinas <-> Montserrat)
# - Refresh the daily case-count figures for several countries
# - Bump the "datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 11:09"

# --- Country name swaps ---
$ws.Range("A128").Value = "Eslovenia"
$ws.Range("A129").Value = "Lituania"

$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Updated statistics ---
# Row 25 - Filipinas
$ws.Range("B25").Value = 147526
$ws.Range("C25").Value = 4002
$ws.Range("D25").Value = 70387
$ws.Range("E25").Value = 74713
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = 2426

# Row 26 - Indonesia
$ws.Range("B26").Value = 132816
$ws.Range("C26").Value = 2098
$ws.Range("D26").Value = 87558
$ws.Range("E26").Value = 39290
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = 5968

# Row 33 - Israel
$ws.Range("B33").Value = 88733
$ws.Range("C33").Value = 582
$ws.Range("D33").Value = 64674
$ws.Range("E33").Value = 23416
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 643

# Row 48 - Polonia
$ws.Range("B48").Value = 54487
$ws.Range("C48").Value = 811
$ws.Range("D48").Value = 37961
$ws.Range("E48").Value = 14682
$ws.Range("G48").Value = 14
$ws.Range("H48").Value = 1844

# Row 71 - Austria
$ws.Range("B71").Value = 22594
$ws.Range("C71").Value = 155
$ws.Range("D71").Value = 20346
$ws.Range("E71").Value = 1523
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 725

# Row 96 - Finlandia
$ws.Range("B96").Value = 7683
$ws.Range("C96").Value = 41
$ws.Range("E96").Value = 300

# Row 111 - Hong Kong
$ws.Range("B111").Value = 4313
$ws.Range("C111").Value = 69
$ws.Range("D111").Value = 3295
$ws.Range("E111").Value = 953
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 65

# Row 123 - Eslovaquia
$ws.Range("B123").Value = 2739
$ws.Range("C123").Value = 49
$ws.Range("D123").Value = 1939
$ws.Range("E123").Value = 769

# Row 128 - (now Eslovenia)
$ws.Range("B128").Value = 2332
$ws.Range("C128").Value = 29
$ws.Range("D128").Value = 1960
$ws.Range("E128").Value = 243
$ws.Range("H128").Value = 129

# Row 129 - (now Lituania)
$ws.Range("B129").Value = 2330
$ws.Range("C129").Value = 21
$ws.Range("D129").Value = 1689
$ws.Range("E129").Value = 560
$ws.Range("H129").Value = 81

# Row 144 - Letonia
$ws.Range("B144").Value = 1307
$ws.Range("C144").Value = 4
$ws.Range("E144").Value = 197

# Row 213 - (now Montserrat)
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 - (now Islas Malvinas)
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
